# Apply the "Error Calculations and Plots" edits to the missing_data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column E (missing-data imputation changes) while row numbers
# --- still match the original layout (these rows are all above the rows
# --- that will be removed below, so row numbers are stable for now).
$ws.Range("E2").Value = -7.2
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E11").Value = -7.9
$ws.Range("E13").ClearContents()
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E25").ClearContents()

# --- Remove the "RM 232" row (row 26) entirely; rows below shift up one.
$ws.Rows(26).Delete()

# --- Remove the "SC 92" row (originally row 28, now row 27 after the
# --- previous deletion) entirely; rows below shift up one again.
$ws.Rows(27).Delete()

# --- Fill in the previously-missing values for the final remaining row
# --- ("SC 232", now row 33).
$ws.Range("C33").Value = 10.4
$ws.Range("E33").Value = -10.7
